$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new blank column before column N ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns.Item(14).Insert()

# Give the newly inserted column (N) a width of 11 characters (raw OOXML width),
# which corresponds to a ColumnWidth property value of 11 - 5/6.
$wsSchedule.Columns.Item(14).ColumnWidth = 10.166666666666666

# --- Make "Repayment schedule" the active sheet/tab and update its selection ---
$wsSchedule.Activate()
$wsSchedule.Range("K17").Select()
